$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows above the current row 5 (E5:G5 "10/30" block),
# shifting the existing "10" and "30" Reynolds rows down to rows 7 and 8.
$ws.Rows.Item(5).Resize(2).Insert()

# New row 5: Re=3
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = -0.0297197
$ws.Range("G5").Value = 4.15251

# New row 6: Re=5
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = -0.0241423
$ws.Range("G6").Value = 2.82355

# Row 7 (previously row 5, Re=10) gains two extra columns H, I
$ws.Range("H7").Value = -0.0159355
$ws.Range("I7").Value = 1.89658

# Append new rows 9 and 10 with further Reynolds numbers
$ws.Range("E9").Value = 50
$ws.Range("F9").Value = 0.00168307
$ws.Range("G9").Value = 1.36542

$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 0.00136985
$ws.Range("G10").Value = 2.16313

# Update the selection to match the author's final active cell
$ws.Range("G6").Select()
